$wb = $excel.ActiveWorkbook

$repoBase = "https://github.com/OpenLocalizationTest/oltest/blob/e7f1a21568cc93c3deeec2430017d0019fa01499"
$hyperlinkColor = 15570276  # OLE BGR encoding of RGB(0x64,0x95,0xED) -> matches the workbook's existing HyperLink font color

function Set-HandoffSheet {
    param($ws, $xlfName, $handoffDatetime)

    # Status: "Handoff transform failed" -> "Ready for handoff"
    $ws.Range("B2").Value = "Ready for handoff"

    # New "Latest Handoff File" entry + hyperlink to the generated xlf
    $ws.Range("C2").Value = $xlfName

    # Latest Handoff Datetime
    $ws.Range("D2").Value = $handoffDatetime

    # Handoff Reason
    $ws.Range("H2").Value = "Include"

    # Rebuild the hyperlinks collection so the new link lands in the middle
    # (between the existing e2e .md link and the .localization-config link),
    # matching the order the report generator produces.
    $ws.Hyperlinks.Delete() | Out-Null
    $ws.Hyperlinks.Add($ws.Range("A2"), "$repoBase/e2e/987b8709-0d39-4f67-9612-2a48fc867e1e.md", "", "", "987b8709-0d39-4f67-9612-2a48fc867e1e.md") | Out-Null
    $ws.Hyperlinks.Add($ws.Range("C2"), "$repoBase/e2e/$xlfName", "", "", $xlfName) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("A3"), "$repoBase/.localization-config", "", "", ".localization-config") | Out-Null

    # Hyperlinks.Add applies Excel's built-in themed "Hyperlink" style; restore
    # the workbook's own underline + blue-font look (same as used elsewhere)
    # for every linked cell.
    $ws.Range("A2").Font.Underline = $true
    $ws.Range("A2").Font.Color = $hyperlinkColor
    $ws.Range("C2").Font.Underline = $true
    $ws.Range("C2").Font.Color = $hyperlinkColor
    $ws.Range("A3").Font.Underline = $true
    $ws.Range("A3").Font.Color = $hyperlinkColor
}

$wsZhCn = $wb.Worksheets.Item("zh-cn")
Set-HandoffSheet $wsZhCn "987b8709-0d39-4f67-9612-2a48fc867e1e.a82daa92b551f26d2f8cad6d66771783ecc88ecc.zh-cn.xlf" "2016-01-17 03:39:53"

$wsDeDe = $wb.Worksheets.Item("de-de")
Set-HandoffSheet $wsDeDe "987b8709-0d39-4f67-9612-2a48fc867e1e.a82daa92b551f26d2f8cad6d66771783ecc88ecc.de-de.xlf" "2016-01-17 03:40:03"

Write-Host "Report regenerated for handoff"
